# Append two new trading-log rows (94, 95) to the bottom of the log sheet,
# mirroring the "TRADING_ATTEMPT" -> "POSITION_OPENED" pair pattern already
# present in the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 94: TRADING_ATTEMPT for TRX
$ws.Range("A94").Value = "2025-10-23T01:29:03.723353"
$ws.Range("B94").Value = "TRADING_ATTEMPT"
$ws.Range("C94").Value = "TRX"
$ws.Range("D94").Value = "UNKNOWN"
$ws.Range("E94").Value = 0.3222621896511093
$ws.Range("K94").Value = "ATTEMPT"
$ws.Range("L94").Value = "Attempting trade 1/1"

# Row 95: resulting POSITION_OPENED for TRX
$ws.Range("A95").Value = "2025-10-23T01:29:05.338613"
$ws.Range("B95").Value = "POSITION_OPENED"
$ws.Range("C95").Value = "TRX"
$ws.Range("D95").Value = "UNKNOWN"
$ws.Range("E95").Value = 0.3222621896511093
$ws.Range("F95").Value = 90
$ws.Range("G95").Value = 1
$ws.Range("H95").Value = 0
$ws.Range("K95").Value = "SUCCESS"
